# Fix: benchmark (UBAH / UBAH Redistribute / Timed Ivy / GTAA) was not
# started 1 year after the first data point, which produced erratic initial
# ROI/MDD/Sharpe figures. This re-pastes the corrected result values that
# came out of re-running the backtest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ROI table (rows 4-7) ---------------------------------------------
# UBAH
$ws.Range("C4").Value2 = 0.83810981027487397
$ws.Range("D4").Value2 = 0.84502819368090598
$ws.Range("E4").Value2 = 0.90925664139523099
$ws.Range("F4").Value2 = 0.97259780002336405
$ws.Range("G4").Value2 = 1.1807117147805
$ws.Range("H4").Value2 = 1.4851402425376501
$ws.Range("H4").Font.Bold = $false

# UBAH Redistribute
$ws.Range("F5").Value2 = 0.98232303679625299
$ws.Range("H5").Value2 = 1.45956731063969

# Timed Ivy
$ws.Range("H6").Value2 = 1.6413053895759999
$ws.Range("H6").Font.Bold = $true

# GTAA
$ws.Range("F7").Value2 = 1.0097627396870299
$ws.Range("H7").Value2 = 1.37289396087158

# --- MDD table (rows 11-14) --------------------------------------------
# UBAH
$ws.Range("C11").Value2 = -0.25750575545606702
$ws.Range("D11").Value2 = -0.25750575545606702
$ws.Range("E11").Value2 = -0.25750575545606702
$ws.Range("F11").Value2 = -0.25750575545606702
$ws.Range("G11").Value2 = -0.25750575545606702
$ws.Range("H11").Value2 = -0.25750575545606702

# Timed Ivy
$ws.Range("C13").Value2 = -0.056664818092968597
$ws.Range("D13").Value2 = -0.066867472204538606
$ws.Range("E13").Value2 = -0.083905675960593895
$ws.Range("F13").Value2 = -0.11790304738742
$ws.Range("G13").Value2 = -0.16947038886165799
$ws.Range("H13").Value2 = -0.16947038886165799

# --- Sharpe Ratio table (rows 18-21) -----------------------------------
# UBAH
$ws.Range("C18").Value2 = 11.593784730166499
$ws.Range("D18").Value2 = 16.160156639056499
$ws.Range("E18").Value2 = 19.678905739952299
$ws.Range("F18").Value2 = 15.6170046920787
$ws.Range("G18").Value2 = 7.9625313937175104
$ws.Range("H18").Value2 = 5.51203251413878

# UBAH Redistribute
$ws.Range("F19").Value2 = 14.912112055510701
$ws.Range("H19").Value2 = 5.7903069057172898

# Timed Ivy
$ws.Range("C20").Value2 = 66.184006361043302
$ws.Range("D20").Value2 = 75.4805911916729
$ws.Range("E20").Value2 = 31.010784506711602
$ws.Range("F20").Value2 = 18.526787433446
$ws.Range("G20").Value2 = 9.9169785230665806
$ws.Range("H20").Value2 = 6.1661769099967598

# GTAA
$ws.Range("D21").Value2 = 36.596582940099601
$ws.Range("E21").Value2 = 23.236473033334899
$ws.Range("F21").Value2 = 23.643242745932302
$ws.Range("H21").Value2 = 8.6241970819911504

# Leave the selection where the author left it when they saved.
$ws.Range("H21").Select()
